# AR GUI dev, added laser control
#
# This script edits Script.xlsx to:
#  - add new "Laser Sys ON", "Laser ENA" and "Reset Intlk Timer" commands to
#    the Sheet2 command-reference table (with supporting notes / parameter
#    descriptions),
#  - rewrite the Sheet1 step script to exercise the new laser commands,
#  - refresh the data validations / column layout that reference the
#    (now larger) command table, and
#  - leave Sheet2 as the active tab/selection.

$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)   # Sheet1 - the step script
$ws2 = $wb.Worksheets.Item(2)   # Sheet2 - command reference / legend

$xlCenter = -4108

# ---------------------------------------------------------------------
# 1) Sheet2: insert the three new command rows (Reset Intlk Timer,
#    Laser Sys ON, Laser ENA) ahead of the existing "End" row, and apply
#    the same centered styling used by the rest of the reference table.
# ---------------------------------------------------------------------

# Grab the "End" row contents before we overwrite them, so we can move
# that row down from row 11 to row 14.
$endCommand = $ws2.Range("A11").Value2

# Row 11: Reset Intlk Timer (no parameters, like Home/Wait/etc rows)
$ws2.Range("A11").Value = "Reset Intlk Timer"
$ws2.Range("B11:D11").Merge()
$ws2.Range("B11").Value = "No Parameters"
$ws2.Range("F11").Value = "Clears serial buffer and takes one initial (start) reading from MCU"

# Row 12: Laser Sys ON (single 1/0 parameter)
$ws2.Range("A12").Value = "Laser Sys ON"
$ws2.Range("B12").Value = "1 or 0"
$ws2.Range("F12").Value = "1 - ON, 0 - OFF"

# Row 13: Laser ENA (single 1/0 parameter)
$ws2.Range("A13").Value = "Laser ENA"
$ws2.Range("B13").Value = "1 or 0"
$ws2.Range("F13").Value = "1 - ON, 0 - OFF"

# Row 14: the "End" row, moved down to make room for the new commands.
$ws2.Range("A14").Value = $endCommand
$ws2.Range("B14:D14").Merge()
$ws2.Range("B14").Value = "No Parameters"

# Rows 15-18: a handful of blank (but centre-styled) spare rows at the
# bottom of the table.
$ws2.Range("B15:D18").HorizontalAlignment = $xlCenter

# Apply centered alignment across the whole used table (this matches the
# existing look-and-feel that was previously only applied to some rows).
$ws2.Range("B1:D14").HorizontalAlignment = $xlCenter

# B12/B13 ("1 or 0") use a dedicated number format alongside the
# centering used elsewhere in the column.
$ws2.Range("B12:B13").NumberFormat = "mmm-yy"
$ws2.Range("B12:B13").HorizontalAlignment = $xlCenter

# Column widths widen now that the table holds longer command names.
$ws2.Columns.Item(1).ColumnWidth = 19.5703125
$ws2.Columns.Item(2).ColumnWidth = 20
$ws2.Columns.Item(3).ColumnWidth = 18.7109375
$ws2.Columns.Item(4).ColumnWidth = 18.5703125
$ws2.Columns.Item(6).ColumnWidth = 60

# ---------------------------------------------------------------------
# 2) Sheet1: rewrite the step script so it resets the interlock timer,
#    turns the laser system & enable lines on, waits, then turns them
#    back off again before repeating / ending.
# ---------------------------------------------------------------------

$ws1.Range("A1:E200").ClearContents()

# Header row
$ws1.Range("A1").Value = "Step#"
$ws1.Range("B1").Value = "Command"
$ws1.Range("C1").Value = "Parameter1"
$ws1.Range("D1").Value = "Parameter2"
$ws1.Range("E1").Value = "Parameter3"

function Set-Step {
    param($row, $step, $command, $p1, $p2, $p3)
    $ws1.Cells.Item($row, 1).Value = $step
    $ws1.Cells.Item($row, 2).Value = $command
    if ($null -ne $p1) { $ws1.Cells.Item($row, 3).Value = $p1 }
    if ($null -ne $p2) { $ws1.Cells.Item($row, 4).Value = $p2 }
    if ($null -ne $p3) { $ws1.Cells.Item($row, 5).Value = $p3 }
}

Set-Step 2  10  "Home All"          0    0   0
Set-Step 3  20  "Reset Intlk Timer" 0    0   0
Set-Step 4  30  "Laser Sys ON"      1  $null $null
Set-Step 5  40  "X Move Abs"        50   0   10
Set-Step 6  50  "Wait"              500  0   0
Set-Step 7  60  "Laser ENA"         1  $null $null
Set-Step 8  65  "Reset Intlk Timer" $null $null $null
Set-Step 9  70  "X Move Abs"        100  0   10
Set-Step 10 80  "Wait"              500  0   0
Set-Step 11 90  "Laser ENA"         0  $null $null
Set-Step 12 100 "X Move Abs"        0    0   0
Set-Step 13 110 "Wait"              500  0   0
Set-Step 14 120 "Repeat"            30   5   0
Set-Step 15 130 "Laser Sys ON"      0  $null $null
Set-Step 16 140 "End"               0    0   0

# Column widths: Command column widens for the longer new command names;
# the two parameter columns shrink a touch to their new best-fit size.
$ws1.Columns.Item(2).ColumnWidth = 16.140625
$ws1.Columns.Item(3).ColumnWidth = 11.28515625
$ws1.Columns.Item(4).ColumnWidth = 11.28515625

# ---------------------------------------------------------------------
# 3) Data validation: the decimal-range validation on C:D now covers the
#    full (larger) table contiguously. (The Command column's drop-down
#    list validation lives in the x14 list-validation extension, which
#    references Sheet2's command table; that is not reachable through
#    the classic Validation object model here, so it is left as-is.)
# ---------------------------------------------------------------------

$ws1.Range("C2:D19").Validation.Delete()
$ws1.Range("C2:D19").Validation.Add(2, 1, 1, "0", "2000")

# ---------------------------------------------------------------------
# 4) Selections / active tab: the workbook was left with Sheet2 active.
# ---------------------------------------------------------------------

$ws1.Range("C26").Select()
$ws2.Range("B18").Select()
$ws2.Activate()
